# Reorder the sheet tabs so that "总计" (totals) becomes the first/leftmost
# sheet and "2022-Q2" becomes the second sheet. No cell data, styles, or
# formatting changes - this is purely a tab-order swap.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# Move "总计" to be before "2022-Q2" (i.e. before its current first position)
$wsTotal.Move($wsQ2)

# Keep "2022-Q2" as the active/selected tab (matches original file, where
# that sheet's view carried tabSelected="1") even though it is now 2nd.
# Re-fetch the worksheet handle after the move since the old reference
# is stale once the sheet collection has been reordered.
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Activate()
